$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its Text storage so Excel does not
# reinterpret values like "510.70" or "1.00" as numbers (losing
# trailing zeros / precision), matching the original inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.305.62'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '3.931.47'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '510.70'
$ws.Range("E5").Value = '  +4.83%  '
$ws.Range("D6").Value = '146.45'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.729'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("E10").Value = '  +3.42%  '
$ws.Range("D11").Value = '0.0000343'
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").Value = '43.08'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '4.564.90'
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").Value = '10.39'
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = '3.925.38'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '14.11'
$ws.Range("E16").Value = '  -2.50%  '
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("E18").Value = '  +6.96%  '
$ws.Range("D19").Value = '19.76'
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("D20").Value = '69.371.52'
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").Value = '432.40'
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("E22").Value = '  -2.60%  '
$ws.Range("D23").Value = '14.53'
$ws.Range("E23").Value = '  -4.53%  '
$ws.Range("D24").Value = '88.29'
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '11.79'
$ws.Range("E25").Value = '  +3.12%  '
$ws.Range("E26").Value = '  +6.14%  '
$ws.Range("D27").Value = '11.08'
$ws.Range("E27").Value = '  -3.53%  '
$ws.Range("D28").Value = '36.74'
$ws.Range("E28").Value = '  -5.20%  '
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").Value = '707.71'
$ws.Range("E30").Value = '  -1.46%  '
$ws.Range("D31").Value = '13.31'
$ws.Range("E31").Value = '  -3.38%  '
$ws.Range("E32").Value = '  -2.47%  '
$ws.Range("D33").Value = '2.87'
$ws.Range("E33").Value = '  -1.68%  '
$ws.Range("D34").Value = '66.61'
$ws.Range("E34").Value = '  +8.63%  '
$ws.Range("D35").Value = '0.445'
$ws.Range("E35").Value = '  +6.37%  '
$ws.Range("D36").Value = '0.0₃0874'
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").Value = '5.96'
$ws.Range("E37").Value = '  -6.56%  '
$ws.Range("D38").Value = '40.46'
$ws.Range("E38").Value = '  -4.60%  '
$ws.Range("D39").Value = '0.149'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("D43").Value = '3.13'
$ws.Range("E43").Value = '  +7.34%  '
$ws.Range("D44").Value = '2.83'
$ws.Range("E44").Value = '  -7.08%  '
$ws.Range("D45").Value = '3.03'
$ws.Range("E45").Value = '  -7.79%  '
$ws.Range("E46").Value = '  +0.75%  '
$ws.Range("D47").Value = '3.33'
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").Value = '0.0₆0355'
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("D49").Value = '3.36'
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("D50").Value = '2.96'
$ws.Range("E50").Value = '  +3.75%  '
$ws.Range("E51").Value = '  -2.23%  '

